# Data refresh for 12/02/2021 snapshot on the "Cycle_2021-2022" sheet (sheet3)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cycle_2021-2022")

# --- Step 1: fill in missing #N/A (=NA()) placeholders for rows 44-58 -------
# (these rows were only partially filled in before; anything still blank in
#  C:E gets an explicit =NA() formula, matching the rest of the sheet's
#  pattern)
$ws.Range("D44:E44").FormulaR1C1   = "=NA()"
$ws.Range("C45:E50").FormulaR1C1   = "=NA()"
$ws.Range("D51:E51").FormulaR1C1   = "=NA()"
$ws.Range("C52:E58").FormulaR1C1   = "=NA()"

# --- Step 2: new "as of" snapshot mini-table header (row 58, columns G:J) ---
$ws.Range("G58").Value = "As of 12/02/2021"
$ws.Range("H58").Value = "n"
$ws.Range("I58").Value = "cGPA"
$ws.Range("J58").Value = "MCAT"

# --- Step 3: append the new 12/02/2021 data rows (59-73) --------------------
# date serial 44532 == 12/2/2021; copy the date style from A58 across the
# whole new block first so every new A-cell gets the same date number format
$ws.Range("A58").Copy()
$ws.Range("A59:A73").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$newRows = @(
    # row, B label,                                   C,     D,     E,     G label,                             H,     I,     J
    @(59, "total applicants",                          1756,  3.69,  511.5, "Total MD/PhD Applicants",           1756,  3.69,  511.5),
    @(60, "withdraw before acceptance (WB)",            37,    3.79,  514.5, "Withdraw Before AC",                 37,    3.79,  514.5),
    @(61, "rejected groups",                            891,   3.67,  511.3, "Defer to Regular MD consideration",  17,    3.81,  515.9),
    @(62, "preliminary rejection",                      $null, $null, $null, "Rejected (PR,PW,RJ)",                891,   3.67,  511.3),
    @(63, "passive withdrawal",                         $null, $null, $null, "At least 1 MD/PhD AC",               160,   3.84,  517.6),
    @(64, "rejection",                                  $null, $null, $null, "Available active (AL,RS,IN)",        479,   3.71,  511.8),
    @(65, "defer to MD app",                            17,    3.81,  515.9, "Available passive (NA,HO)",          172,   3.57,  504.1),
    @(66, "at least 1 MD/PhD acceptance",                160,   3.84,  517.6, $null,                                $null, $null, $null),
    @(67, "available active",                           479,   3.71,  511.8, $null,                                $null, $null, $null),
    @(68, "request secondary",                          $null, $null, $null, $null,                                $null, $null, $null),
    @(69, "interview scheduled",                        $null, $null, $null, $null,                                $null, $null, $null),
    @(70, "available passive",                          172,   3.57,  504.1, $null,                                $null, $null, $null),
    @(71, "no action",                                  $null, $null, $null, $null,                                $null, $null, $null),
    @(72, "hold",                                       $null, $null, $null, $null,                                $null, $null, $null),
    @(73, "available",                                  $null, $null, $null, $null,                                $null, $null, $null)
)

foreach ($r in $newRows) {
    $row = $r[0]

    $ws.Cells.Item($row, 1).Value = 44532
    $ws.Cells.Item($row, 2).Value = $r[1]

    if ($null -eq $r[2]) {
        $ws.Range("C" + $row + ":E" + $row).FormulaR1C1 = "=NA()"
    } else {
        $ws.Cells.Item($row, 3).Value = $r[2]
        $ws.Cells.Item($row, 4).Value = $r[3]
        $ws.Cells.Item($row, 5).Value = $r[4]
    }

    if ($null -ne $r[5]) {
        $ws.Cells.Item($row, 7).Value  = $r[5]
        $ws.Cells.Item($row, 8).Value  = $r[6]
        $ws.Cells.Item($row, 9).Value  = $r[7]
        $ws.Cells.Item($row, 10).Value = $r[8]
    }
}

# --- Step 4: restore the view/selection state recorded in the saved file ----
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("K51").Select()
